# Apply cryptocurrency price/volume updates to Sheet1 (generated from commit diff)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new values are numeric-looking strings (e.g. "0.999", "119.20").
# Excel would otherwise auto-convert these to numbers on assignment (losing trailing zeros,
# e.g. "119.20" -> 119.2), so force plain-text storage on just these cells first, matching
# the source workbook where every Price/Volume cell is stored as text.
$textCells = @("D5", "D6", "D14", "D20", "D21", "D22", "D24", "D25", "D27", "D31", "D32", "D33", "D34", "D35", "D37", "D38", "D43", "D46")
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range("D2").Value = "57.069.40"
$ws.Range("E2").Value = "  -5.67%  "
$ws.Range("D3").Value = "2.892.25"
$ws.Range("E3").Value = "  -3.27%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "546.29"
$ws.Range("E5").Value = "  -3.08%  "
$ws.Range("D6").Value = "124.25"
$ws.Range("E6").Value = "  -1.15%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("E8").Value = "  +0.94%  "
$ws.Range("D9").Value = "2.885.26"
$ws.Range("E9").Value = "  -3.35%  "
$ws.Range("E10").Value = "  -8.51%  "
$ws.Range("E11").Value = "  -8.63%  "
$ws.Range("E12").Value = "  -1.02%  "
$ws.Range("E13").Value = "  -6.58%  "
$ws.Range("D14").Value = "31.99"
$ws.Range("E14").Value = "  -1.80%  "
$ws.Range("E15").Value = "  +1.16%  "
$ws.Range("D16").Value = "3.363.76"
$ws.Range("E16").Value = "  -3.29%  "
$ws.Range("D17").Value = "2.887.99"
$ws.Range("E17").Value = "  -3.13%  "
$ws.Range("E18").Value = "  +5.40%  "
$ws.Range("D19").Value = "57.017.90"
$ws.Range("E19").Value = "  -5.91%  "
$ws.Range("D20").Value = "401.58"
$ws.Range("E20").Value = "  -6.54%  "
$ws.Range("D21").Value = "12.71"
$ws.Range("E21").Value = "  -2.25%  "
$ws.Range("D22").Value = "0.666"
$ws.Range("E22").Value = "  +1.14%  "
$ws.Range("E23").Value = "  -4.56%  "
$ws.Range("D24").Value = "12.57"
$ws.Range("E24").Value = "  -2.96%  "
$ws.Range("D25").Value = "77.26"
$ws.Range("E25").Value = "  -1.50%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  +0.11%  "
$ws.Range("E28").Value = "  -1.34%  "
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("E30").Value = "  +1.19%  "
$ws.Range("D31").Value = "24.49"
$ws.Range("E31").Value = "  -2.89%  "
$ws.Range("D32").Value = "5.90"
$ws.Range("E32").Value = "  -1.93%  "
$ws.Range("D33").Value = "0.0976"
$ws.Range("E33").Value = "  +4.56%  "
$ws.Range("D34").Value = "0.911"
$ws.Range("E34").Value = "  -4.10%  "
$ws.Range("D35").Value = "5.40"
$ws.Range("E35").Value = "  -2.11%  "
$ws.Range("E36").Value = "  -11.75%  "
$ws.Range("D37").Value = "47.98"
$ws.Range("E37").Value = "  -2.77%  "
$ws.Range("D38").Value = "8.19"
$ws.Range("E38").Value = "  +5.31%  "
$ws.Range("E39").Value = "  -5.80%  "
$ws.Range("E40").Value = "  -1.20%  "
$ws.Range("E41").Value = "  -6.42%  "
$ws.Range("D42").Value = "2.614.20"
$ws.Range("E42").Value = "  -1.73%  "
$ws.Range("D43").Value = "360.07"
$ws.Range("E43").Value = "  -3.86%  "
$ws.Range("E44").Value = "  -0.91%  "
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("D46").Value = "119.20"
$ws.Range("E46").Value = "  -0.08%  "
$ws.Range("E47").Value = "  +0.89%  "
$ws.Range("E48").Value = "  -3.31%  "
$ws.Range("E49").Value = "  -2.33%  "
$ws.Range("E50").Value = "  -3.95%  "
$ws.Range("E51").Value = "  -3.27%  "
